# Add a new "30-jun" data column (Q) to the weekly tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added date column.
$ws.Range("Q1").Value = "30-jun"

# New data values for the added column, one per product row (rows 2-11).
$values = @(16, 15, 8, 11, 13, 21, 10, 17, 7, 15)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Range("Q$row")
    # Match the formatting used by the rest of the numeric columns
    # (centered, integer number format) before writing the value.
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
    $cell.Value = $values[$i]
}

# Restore the scroll position / active selection recorded for the sheet.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P19").Select()
